$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content fix: GegnerTests numbering correction ---
# Row 11 (Anforderung "Der Gegner muss mit Spieler kollidieren können") was mislabeled
# "2.2.3.2.6"; correct numbering is "2.2.3.5.2" (it belongs in the 2.2.3.5.x collision
# group, right after 2.2.3.5.1).
$ws.Range("A11").Value = "2.2.3.5.2"

# --- Style cleanup: B6 carried a redundant "applyFill" flag that none of the other
# description cells (B2:B5, B7:B11) have. Re-apply the same formatting as a sibling
# cell so it shares their (cleaner) style record. ---
$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Remove the extra trailing blank placeholder row (sheet previously reserved
# rows 12-16 with alternating styles; only 12-15 are actually needed/used). ---
[void]$ws.Rows.Item(16).Delete()

# --- Selection bookkeeping to match the saved workbook's last cursor position. ---
[void]$ws.Range("B14").Select()
